$wb = $excel.ActiveWorkbook
$daily = $wb.Worksheets.Item("Daily")
$cum = $wb.Worksheets.Item("Cum")

# Add new row 15 to Daily sheet
$daily.Range("A15").Value = 44540
$daily.Range("B15").Value = 443
$daily.Range("C15").Value = 1
$daily.Range("D15").Value = 4
$daily.Range("F15").Formula = "=SUM(B15:E15)"

# Add new row 15 to Cum sheet
$cum.Range("A15").Value = 44540
$cum.Range("B15").Formula = "=SUM(B14+Daily!B15)"
$cum.Range("C15").Formula = "=SUM(C14+Daily!C15)"
$cum.Range("D15").Formula = "=SUM(D14+Daily!D15)"
$cum.Range("E15").Formula = "=SUM(E14+Daily!E15)"
$cum.Range("F15").Formula = "=SUM(B15:E15)"

$daily.Select()
$daily.Range("W45").Select()

$cum.Select()
$cum.Range("A22").Select()
$cum.Range("W72").Select()
